$wb = $excel.ActiveWorkbook

# --- Sheet1 (Hardware) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A3:B8").ClearContents()
$ws1.Range("A2").Value = "H501#G103"

# --- Sheet2 (Options) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A45:D56").ClearContents()

$ws2.Range("A2").Value = "J580"
$ws2.Range("A3").Value = "J749"
$ws2.Range("A4").Value = "S728"
$ws2.Range("A5").Value = "S985"
$ws2.Range("A6").Value = "J801"
$ws2.Range("A7").Value = "J802#11"
$ws2.Range("A8").Value = "J803"
$ws2.Range("A9").Value = "J804"
$ws2.Range("A10").Value = "J807"
$ws2.Range("A11").Value = "J819"
$ws2.Range("A12").Value = "J828"
$ws2.Range("A13").Value = "J829"
$ws2.Range("A14").Value = "J830"
$ws2.Range("A15").Value = "J835"
$ws2.Range("A16").Value = "J836"
$ws2.Range("A17").Value = "J838"
$ws2.Range("A18").Value = "J841"
$ws2.Range("A19").Value = "J842"
$ws2.Range("A20").Value = "J846"
$ws2.Range("A21").Value = "J848"
$ws2.Range("A22").Value = "J850"
$ws2.Range("A23").Value = "J853"
$ws2.Range("A24").Value = "J854"
$ws2.Range("A25").Value = "J873"
$ws2.Range("A26").Value = "J876"
$ws2.Range("A27").Value = "J882"
$ws2.Range("A28").Value = "J887"
$ws2.Range("A29").Value = "J890"
$ws2.Range("A30").Value = "J893"
$ws2.Range("A31").Value = "J894"
$ws2.Range("A32").Value = "J895"
$ws2.Range("A33").Value = "J900"
$ws2.Range("A34").Value = "J913"
$ws2.Range("A35").Value = "J917"
$ws2.Range("A36").Value = "J930"
$ws2.Range("A37").Value = "J937"
$ws2.Range("A38").Value = "J948"
$ws2.Range("A39").Value = "v953"
$ws2.Range("A40").Value = "J956"
$ws2.Range("A41").Value = "J965"
$ws2.Range("A42").Value = "J971"
$ws2.Range("A43").Value = "J981"
$ws2.Range("A44").Value = "R955"
